$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.472.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "'1.795.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'223.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").Value = "'0.550"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'32.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("D9").Value = "'0.287"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.07%  "

$ws.Range("D10").Value = "'0.0703"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.42%  "

$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("D12").Value = "'2.055.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.64%  "

$ws.Range("D13").Value = "'10.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.82%  "

$ws.Range("D14").Value = "'1.785.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.02%  "

$ws.Range("D15").Value = "'0.637"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.44%  "

$ws.Range("D16").Value = "'34.510.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").Value = "'4.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").Value = "'68.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.95%  "

$ws.Range("D19").Value = "'249.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.42%  "

$ws.Range("D20").Value = "'0.0₃0794"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.82%  "

$ws.Range("D21").Value = "'11.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.92%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "'4.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.06%  "

$ws.Range("D24").Value = "'2.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "'161.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.48%  "

$ws.Range("D26").Value = "'16.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.75%  "

$ws.Range("D27").Value = "'7.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.68%  "

$ws.Range("E28").Value = "  -0.62%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").Value = "'567.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +993.17%  "

$ws.Range("D31").Value = "'0.0522"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("D32").Value = "'3.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.64%  "

$ws.Range("D33").Value = "'1.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.75%  "

$ws.Range("D34").Value = "'3.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("D35").Value = "'1.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.54%  "

$ws.Range("D36").Value = "'1.424.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.87%  "

$ws.Range("D37").Value = "'1.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("D38").Value = "'0.634"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("D39").Value = "'0.0189"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.69%  "

$ws.Range("D40").Value = "'83.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("D41").Value = "'0.948"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.22%  "

$ws.Range("D42").Value = "'2.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.19%  "

$ws.Range("D43").Value = "'2.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("D44").Value = "'2.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.38%  "

$ws.Range("D45").Value = "'6.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.12%  "

$ws.Range("E46").Value = "  -0.99%  "

$ws.Range("D47").Value = "'0.0497"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.35%  "

$ws.Range("D48").Value = "'1.946.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'105.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.41%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'12.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("E51").Value = "  -0.03%  "
